# Sample Project / Main.xlsx - "Rules" sheet, row 11 (the R40 rule row).
# B11 changes from the text "R40" to the text "1" - still a text value
# (not a number), so force the text/General quoting via the leading
# apostrophe the same way a user typing into the Excel UI would, rather
# than letting automatic number recognition turn it into the numeral 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
